$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the kernel SVR parameters (columns K, L, M)
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Their corresponding values on row 2
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Match the selection left behind in the saved workbook
$ws.Range("K8").Select()
